$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.889.41'
$ws.Range("E2").Value = '  -0.22%  '

$ws.Range("D3").Value = '1.639.06'
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.04'
$ws.Range("E5").Value = '  -0.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5028'
$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  -0.18%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2571'
$ws.Range("E8").Value = '  -0.31%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06390'
$ws.Range("E9").Value = '  -0.67%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.53'
$ws.Range("E10").Value = '  -0.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07799'
$ws.Range("E11").Value = '  +0.79%  '

$ws.Range("D12").Value = '1.652.30'
$ws.Range("E12").Value = '  +0.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.270'
$ws.Range("E13").Value = '  +0.51%  '

$ws.Range("D14").Value = '1.865.46'
$ws.Range("E14").Value = '  +0.28%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5420'
$ws.Range("E15").Value = '  -0.64%  '

$ws.Range("D16").Value = '0.0₅7866'
$ws.Range("E16").Value = '  -0.99%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.69'
$ws.Range("E17").Value = '  +1.68%  '

$ws.Range("D18").Value = '25.950.29'
$ws.Range("E18").Value = '  +0.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  -0.27%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '197.80'
$ws.Range("E20").Value = '  -3.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.383'
$ws.Range("E21").Value = '  +1.70%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.931'
$ws.Range("E22").Value = '  -0.65%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.968'
$ws.Range("E23").Value = '  +0.37%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.007'
$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.878'
$ws.Range("E25").Value = '  -2.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '139.91'
$ws.Range("E26").Value = '  -0.98%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1142'
$ws.Range("E27").Value = '  -1.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.843'
$ws.Range("E28").Value = '  +1.15%  '

$ws.Range("E29").Value = '  -0.66%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.242'
$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04875'
$ws.Range("E31").Value = '  -3.68%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.262'
$ws.Range("E32").Value = '  -0.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.193'
$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.532'
$ws.Range("E34").Value = '  -0.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.371'
$ws.Range("E35").Value = '  +1.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8895'
$ws.Range("E36").Value = '  -0.61%  '

$ws.Range("E37").Value = '  -0.48%  '

$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.134.00'
$ws.Range("E38").Value = '  -0.36%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5538'
$ws.Range("E39").Value = '  -1.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01562'
$ws.Range("E40").Value = '  +0.09%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.007'
$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.691'
$ws.Range("E42").Value = '  +0.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8182'
$ws.Range("E43").Value = '  +0.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.38'
$ws.Range("E44").Value = '  +0.15%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.0₈123'
$ws.Range("E45").Value = '  +10.76%  '

$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.776.13'
$ws.Range("E46").Value = '  +0.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4521'
$ws.Range("E47").Value = '  -0.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.009'
$ws.Range("E48").Value = '  +0.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.18'
$ws.Range("E49").Value = '  +0.72%  '

$ws.Range("E50").Value = '  +1.34%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.007'
$ws.Range("E51").Value = '  -0.03%  '
